$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(12,75),
    @(12,311),
    @(12,345),
    @(12,14),
    @(12,130),
    @(12,2),
    @(12,124),
    @(12,78),
    @(12,258),
    @(12,136),
    @(12,105),
    @(12,121),
    @(12,323),
    @(12,357),
    @(12,15),
    @(12,211),
    @(12,147),
    @(12,3),
    @(13,156),
    @(13,255),
    @(13,321),
    @(13,296),
    @(13,286),
    @(13,228),
    @(13,152),
    @(13,340),
    @(13,307),
    @(13,324),
    @(13,159),
    @(13,92),
    @(13,49),
    @(13,325),
    @(13,252),
    @(13,349),
    @(13,98),
    @(13,316),
    @(14,238),
    @(14,273),
    @(14,23),
    @(14,170),
    @(14,194),
    @(14,6),
    @(14,155),
    @(14,270),
    @(14,330),
    @(14,175),
    @(14,226),
    @(14,101),
    @(14,86),
    @(14,344),
    @(14,84),
    @(14,268),
    @(14,79),
    @(14,308),
    @(15,281),
    @(15,343),
    @(15,302),
    @(15,334),
    @(15,317),
    @(15,91),
    @(15,193),
    @(15,209),
    @(15,48),
    @(15,143),
    @(15,54),
    @(15,149),
    @(15,201),
    @(15,164),
    @(15,333),
    @(15,360),
    @(15,197),
    @(15,104),
    @(16,181),
    @(16,7),
    @(16,350),
    @(16,306),
    @(16,1),
    @(16,240),
    @(16,46),
    @(16,284),
    @(16,150),
    @(16,294),
    @(16,16),
    @(16,56),
    @(16,26),
    @(16,244),
    @(16,331),
    @(16,30),
    @(16,293),
    @(16,245)
)

$startRow = 199
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

$ws.Range("B289").Select() | Out-Null
